$wb = $excel.ActiveWorkbook

# --- Update the "Date" metadata value on the Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2022-03-09T14:32:46-05:00"

# --- Remove the duplicate "HRESCH / healthcare research" row (old row 11)
#     from the "Include from ActReason" sheet; rows 12-13 shift up to 11-12 ---
$inc = $wb.Worksheets.Item("Include from ActReason")
$inc.Range("A11").EntireRow.Delete()
